$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 3492.1667
$ws.Range("I86").Value = 2170.6
$ws.Range("J86").Value = 10100
$ws.Range("K86").Value = 2170.6
$ws.Range("L86").Value = 10100
$ws.Range("M86").Value = -1047.6
$ws.Range("N86").Value = -12346

# Row 89
$ws.Range("H89").Value = 3492.1667
$ws.Range("I89").Value = 2170.6
$ws.Range("J89").Value = 10100
$ws.Range("K89").Value = 10853
$ws.Range("L89").Value = 50500
$ws.Range("M89").Value = -5237
$ws.Range("N89").Value = -61732

# Row 112
$ws.Range("H112").Value = 20001436
$ws.Range("I112").Value = 250000900
$ws.Range("J112").Value = 1483
$ws.Range("K112").Value = 750002700
$ws.Range("L112").Value = 4449
$ws.Range("M112").Value = -750001592
$ws.Range("N112").Value = -6665

# Row 129
$ws.Range("H129").Value = 1051.9246
$ws.Range("I129").Value = 400
$ws.Range("J129").Value = 1105.1428
$ws.Range("K129").Value = 1200
$ws.Range("L129").Value = 3315.4284
$ws.Range("M129").Value = 3800
$ws.Range("N129").Value = -13315.4284

# Row 138
$ws.Range("H138").Value = 4315.026
$ws.Range("I138").Value = 2156.4119
$ws.Range("J138").Value = 4926.6333
$ws.Range("K138").Value = 6469.2357
$ws.Range("L138").Value = 14779.8999
$ws.Range("M138").Value = -1329.2357
$ws.Range("N138").Value = -25059.8999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 25975.75
$ws.Range("I32").Value = 44104.8
$ws.Range("J32").Value = 22618.518
$ws.Range("K32").Value = 44104.8
$ws.Range("L32").Value = 22618.518
$ws.Range("M32").Value = -43817.8
$ws.Range("N32").Value = -23192.518

# Row 102
$ws.Range("H102").Value = 2167.25
$ws.Range("I102").Value = 2001.3334
$ws.Range("J102").Value = 2333.1667
$ws.Range("K102").Value = 2001.3334
$ws.Range("L102").Value = 2333.1667
$ws.Range("M102").Value = -379.3334
$ws.Range("N102").Value = -5577.1667

# Row 132
$ws.Range("H132").Value = 2708.257
$ws.Range("I132").Value = 1120.7084
$ws.Range("J132").Value = 6172
$ws.Range("K132").Value = 3362.1252
$ws.Range("L132").Value = 18516
$ws.Range("M132").Value = -832.1251999999999
$ws.Range("N132").Value = -23576

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 1839056.6
$ws.Range("I7").Value = 6677667.5
$ws.Range("J7").Value = 24577.5
$ws.Range("K7").Value = 6677667.5
$ws.Range("L7").Value = 24577.5
$ws.Range("M7").Value = -6677554.5
$ws.Range("N7").Value = -24803.5

# Row 94
$ws.Range("H94").Value = 11629154
$ws.Range("I94").Value = 13158964
$ws.Range("J94").Value = 2596
$ws.Range("K94").Value = 13158964
$ws.Range("L94").Value = 2596
$ws.Range("M94").Value = -13158513
$ws.Range("N94").Value = -3498

# Row 107
$ws.Range("H107").Value = 1375.2142
$ws.Range("I107").Value = 1204.8182
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1204.8182
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 715.1818000000001
$ws.Range("N107").Value = -5840

# Row 134
$ws.Range("H134").Value = 2560.137
$ws.Range("I134").Value = 1686.75
$ws.Range("J134").Value = 4722.8096
$ws.Range("K134").Value = 5060.25
$ws.Range("L134").Value = 14168.4288
$ws.Range("M134").Value = -2525.25
$ws.Range("N134").Value = -19238.4288

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7299.9287
$ws.Range("I31").Value = 1800
$ws.Range("J31").Value = 8216.583000000001
$ws.Range("K31").Value = 1800
$ws.Range("L31").Value = 8216.583000000001
$ws.Range("M31").Value = -1505
$ws.Range("N31").Value = -8806.583000000001

# Row 34
$ws.Range("H34").Value = 7299.9287
$ws.Range("I34").Value = 1800
$ws.Range("J34").Value = 8216.583000000001
$ws.Range("K34").Value = 1800
$ws.Range("L34").Value = 8216.583000000001
$ws.Range("M34").Value = -1598
$ws.Range("N34").Value = -8620.583000000001

# Row 58
$ws.Range("H58").Value = 2104.4546
$ws.Range("I58").Value = 1601.3693
$ws.Range("J58").Value = 4829.5
$ws.Range("K58").Value = 1601.3693
$ws.Range("L58").Value = 4829.5
$ws.Range("M58").Value = -1398.3693
$ws.Range("N58").Value = -5235.5

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 99
$ws.Range("H99").Value = 3067.5356
$ws.Range("I99").Value = 1736.8948
$ws.Range("J99").Value = 5876.6665
$ws.Range("K99").Value = 1736.8948
$ws.Range("L99").Value = 5876.6665
$ws.Range("M99").Value = -238.8948
$ws.Range("N99").Value = -8872.666499999999

# Row 126
$ws.Range("H126").Value = 3067.5356
$ws.Range("I126").Value = 1736.8948
$ws.Range("J126").Value = 5876.6665
$ws.Range("K126").Value = 5210.6844
$ws.Range("L126").Value = 17629.9995
$ws.Range("M126").Value = -2740.6844
$ws.Range("N126").Value = -22569.9995

# Row 132
$ws.Range("H132").Value = 2271.2622
$ws.Range("I132").Value = 1780.3265
$ws.Range("J132").Value = 4275.9165
$ws.Range("K132").Value = 5340.979499999999
$ws.Range("L132").Value = 12827.7495
$ws.Range("M132").Value = -2810.979499999999
$ws.Range("N132").Value = -17887.7495

# Row 136
$ws.Range("H136").Value = 2104.4546
$ws.Range("I136").Value = 1601.3693
$ws.Range("J136").Value = 4829.5
$ws.Range("K136").Value = 4804.1079
$ws.Range("L136").Value = 14488.5
$ws.Range("M136").Value = -2254.1079
$ws.Range("N136").Value = -19588.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1967.4
$ws.Range("I5").Value = 1263.25
$ws.Range("J5").Value = 2298.7646
$ws.Range("K5").Value = 3789.75
$ws.Range("L5").Value = 6896.293799999999
$ws.Range("M5").Value = -3677.75
$ws.Range("N5").Value = -7120.293799999999

# Row 128
$ws.Range("H128").Value = 123333.336
$ws.Range("I128").Value = 123333.336
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 370000.008
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -365020.008

# Row 135
$ws.Range("H135").Value = 1967.4
$ws.Range("I135").Value = 1263.25
$ws.Range("J135").Value = 2298.7646
$ws.Range("K135").Value = 11369.25
$ws.Range("L135").Value = 20688.8814
$ws.Range("M135").Value = -8834.25
$ws.Range("N135").Value = -25758.8814

$ws = $wb.Worksheets.Item("GSM")
# Row 134
$ws.Range("H134").Value = 62820.2
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 62820.2
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 188460.6
$ws.Range("N134").Value = -193530.6

# Row 136
$ws.Range("H136").Value = 46663
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 46663
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 139989
$ws.Range("N136").Value = -145089

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3448.5483
$ws.Range("I122").Value = 2737.2222
$ws.Range("J122").Value = 8250
$ws.Range("K122").Value = 8211.6666
$ws.Range("L122").Value = 24750
$ws.Range("M122").Value = -5761.6666
$ws.Range("N122").Value = -29650

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 146250
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 146250
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 146250
$ws.Range("N46").Value = -146712

# Row 133
$ws.Range("H133").Value = 43388.8
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 43388.8
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 43388.8
$ws.Range("N133").Value = -53508.8

# Row 134
$ws.Range("H134").Value = 146250
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 146250
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 438750
$ws.Range("N134").Value = -443820

# Row 136
$ws.Range("H136").Value = 2879.5386
$ws.Range("I136").Value = 1150.4286
$ws.Range("J136").Value = 7280.909
$ws.Range("K136").Value = 3451.2858
$ws.Range("L136").Value = 21842.727
$ws.Range("M136").Value = -901.2857999999997
$ws.Range("N136").Value = -26942.727
